# Updated cryptos list on Tue Jun 13 20:23:27 UTC 2023 with GitHub Actions
#
# Refreshes the scraped Price (column D) / Volume(1h) (column E) figures
# for each coin row on the sheet.
#
# Rows 12/13 additionally swap rank position between Solana and Polygon in
# this refresh (Solana moved ahead of Polygon), so the Coin name, Link and
# Price/Volume values for those two rows are fully rewritten rather than
# just the numbers.
#
# Column D holds plain scraped text (e.g. "25.873.47", "237.41", "1.001")
# rather than real numbers, and some of the new values (e.g. "237.57",
# "1.000", "14.92") would otherwise be auto-converted to floating point
# numbers by Excel on assignment. Prefixing with an apostrophe - exactly
# what typing '237.57 into a cell does - keeps them stored as literal
# Text, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$q = "'"

$ws.Range("D2").Value = $q + '25.869.23'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = $q + '1.739.18'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = $q + '237.57'
$ws.Range("E5").Value = '  +3.72%  '
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").Value = $q + '0.5172'
$ws.Range("E7").Value = '  -1.61%  '
$ws.Range("D8").Value = $q + '0.2747'
$ws.Range("E8").Value = '  +0.31%  '
$ws.Range("D9").Value = $q + '0.06156'
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").Value = $q + '1.743.33'
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("D11").Value = $q + '0.07185'
$ws.Range("E11").Value = '  +1.37%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = $q + '14.92'
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = $q + '0.6414'
$ws.Range("E13").Value = '  +0.68%  '
$ws.Range("E14").Value = '  +1.70%  '
$ws.Range("D15").Value = $q + '77.42'
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").Value = $q + '1.000'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").Value = $q + '25.885.34'
$ws.Range("E18").Value = '  +0.26%  '
$ws.Range("D19").Value = $q + '11.69'
$ws.Range("E19").Value = '  +1.62%  '
$ws.Range("D20").Value = $q + '0.000006765'
$ws.Range("E20").Value = '  +1.94%  '
$ws.Range("D21").Value = $q + '1.967.87'
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").Value = $q + '4.277'
$ws.Range("E22").Value = '  +0.94%  '
$ws.Range("D23").Value = $q + '8.623'
$ws.Range("E23").Value = '  -1.88%  '
$ws.Range("D24").Value = $q + '5.261'
$ws.Range("E24").Value = '  +2.06%  '
$ws.Range("D25").Value = $q + '139.01'
$ws.Range("E25").Value = '  -0.88%  '
$ws.Range("D26").Value = $q + '1.514'
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").Value = $q + '15.13'
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("D28").Value = $q + '1.755'
$ws.Range("E28").Value = '  -1.31%  '
$ws.Range("D29").Value = $q + '105.61'
$ws.Range("E29").Value = '  +3.47%  '
$ws.Range("D30").Value = $q + '3.908'
$ws.Range("E30").Value = '  +5.39%  '
$ws.Range("D31").Value = $q + '0.08267'
$ws.Range("E31").Value = '  -0.60%  '
$ws.Range("D32").Value = $q + '3.690'
$ws.Range("E32").Value = '  +4.52%  '
$ws.Range("D33").Value = $q + '0.04597'
$ws.Range("E33").Value = '  +2.98%  '
$ws.Range("D34").Value = $q + '2.644'
$ws.Range("E34").Value = '  +1.22%  '
$ws.Range("D35").Value = $q + '0.9872'
$ws.Range("E35").Value = '  +1.54%  '
$ws.Range("D36").Value = $q + '0.6165'
$ws.Range("E36").Value = '  -0.29%  '
$ws.Range("D37").Value = $q + '2.678'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").Value = $q + '0.01602'
$ws.Range("E38").Value = '  +1.99%  '
$ws.Range("D39").Value = $q + '1.921'
$ws.Range("E39").Value = '  +1.33%  '
$ws.Range("D40").Value = $q + '1.000'
$ws.Range("E40").Value = '  +0.12%  '
$ws.Range("D41").Value = $q + '97.65'
$ws.Range("E41").Value = '  -2.25%  '
$ws.Range("D42").Value = $q + '0.3832'
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = $q + '0.7390'
$ws.Range("E43").Value = '  +2.18%  '
$ws.Range("D44").Value = $q + '4.980'
$ws.Range("E44").Value = '  -0.88%  '
$ws.Range("D45").Value = $q + '0.1129'
$ws.Range("E45").Value = '  +0.84%  '
$ws.Range("D46").Value = $q + '6.216'
$ws.Range("E46").Value = '  +0.44%  '
$ws.Range("D47").Value = $q + '0.05241'
$ws.Range("E47").Value = '  -1.56%  '
$ws.Range("D48").Value = $q + '54.70'
$ws.Range("E48").Value = '  +2.58%  '
$ws.Range("D49").Value = $q + '30.38'
$ws.Range("E49").Value = '  +1.43%  '
$ws.Range("D50").Value = $q + '7.552'
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("D51").Value = $q + '0.3399'
$ws.Range("E51").Value = '  +0.27%  '
